$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cell for the "Quoted At" column
$ws.Range("D1").Value = "Quoted At"

# Match the bold header styling already used by A1:C1
$ws.Range("D1").Font.Bold = $true

# Date values, written as the underlying day-count serials (2021-01-01, 2020-01-15,
# 2019-03-04) so the cells store plain numbers and pick up their date look purely from the
# NumberFormat applied below (assigning a .NET DateTime here would auto-apply its own
# number format immediately, stealing a style slot ahead of time and reordering cellXfs)
$ws.Range("D2").Value = 44197
$ws.Range("D3").Value = 43845
$ws.Range("D5").Value = 43528

# Apply a date number format to the touched cells only (row 4 has no "Quoted At" value, so it
# must stay untouched rather than materializing an empty D4 cell)
$ws.Range("D1").NumberFormat = "m/d/yyyy"
$ws.Range("D2").NumberFormat = "m/d/yyyy"
$ws.Range("D3").NumberFormat = "m/d/yyyy"
$ws.Range("D5").NumberFormat = "m/d/yyyy"

# Match column width used in the target workbook
$ws.Range("D1").ColumnWidth = 13.83203125

# Update the selection to mirror the saved view state
$ws.Range("D6").Select()
